$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append (columns: A=datetime serial, B=symbol, C=open, D=high, E=low, F=close, G=volume)
$newRows = @(
    @(45107, "ECONOMICS:CNCBBS", 41806284000000, 41806284000000, 41806284000000, 41806284000000, 0),
    @(45138, "ECONOMICS:CNCBBS", 40809168000000, 40809168000000, 40809168000000, 40809168000000, 0),
    @(45169, "ECONOMICS:CNCBBS", 41684045000000, 41684045000000, 41684045000000, 41684045000000, 0),
    @(45199, "ECONOMICS:CNCBBS", 42735489000000, 42735489000000, 42735489000000, 42735489000000, 0)
)

$startRow = 303
$ws.Range("A302").Copy()
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Range("A$r").PasteSpecial(-4122)
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}
$excel.CutCopyMode = 0
